# Caracteristicas_Dominio07.xlsx - "atualizacoes dos scripts do nanda"
#
# Acrescenta seis novas características (linhas) à tabela de
# "Interação social prejudicada" na planilha Plan1, seguindo o mesmo
# padrão das linhas já existentes (colunas A/C/E fixas, B = código do
# diagnóstico, D = texto da característica).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# A última linha já preenchida na planilha é a 331; as novas entradas
# continuam a partir da linha 332.
$firstNewRow = 332

$newItems = @(
    "Comportamentos de interação social malsucedidos ",
    "Desconforto em situações sociais ",
    "Incapacidade de comunicar uma sensação satisfatória de envolvimento social (p. ex., pertinência, cuidado, interesse, história compartilhada) ",
    "Incapacidade de receber uma sensação satisfatória de envolvimento social (p. ex., pertinência, cuidado, interesse, história compartilhada) ",
    "Interação disfuncional com outras pessoas ",
    "Relato familiar de mudança na interação (p. ex., estilo, padrão) "
)

$row = $firstNewRow
foreach ($item in $newItems) {
    $ws.Cells.Item($row, 1).Value = "###"
    $ws.Cells.Item($row, 2).Value = 115
    $ws.Cells.Item($row, 3).Value = "%%%"
    $ws.Cells.Item($row, 4).Value = $item
    $ws.Cells.Item($row, 5).Value = "$$$"
    $row = $row + 1
}

$lastNewRow = $row - 1

# Reflete a seleção/posição de tela deixada pelo autor após digitar as
# novas linhas (coluna A, da última linha antiga até a última nova).
$ws.Activate()
$ws.Range("A" + ($firstNewRow - 1) + ":A" + $lastNewRow).Select()
$excel.ActiveWindow.ScrollRow = 313
